$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1 - copy formatting from the
# neighboring header cell (G1) so it matches the existing bold/border style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill the new Save column (H2:H6) with the value 1 for every data row.
$ws.Range("H2:H6").Value = 1
